# Update the "Jogos da Semana" odds sheet:
#   * refresh several betting-odds values on existing rows 3,4,5,6,9,10
#   * append a new fixture as row 11 (New Mexico vs Las Vegas Lights)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values (rows 3-10) ---
$ws.Cells.Item(3, 7).Value = 1.57
$ws.Cells.Item(3, 8).Value = 3.8
$ws.Cells.Item(3, 9).Value = 6.25
$ws.Cells.Item(3, 21).Value = 2.1
$ws.Cells.Item(3, 22).Value = 1.67
$ws.Cells.Item(3, 25).Value = 8.5

$ws.Cells.Item(4, 13).Value = 1.07
$ws.Cells.Item(4, 15).Value = 1.41
$ws.Cells.Item(4, 16).Value = 2.62

$ws.Cells.Item(5, 7).Value = 1.45
$ws.Cells.Item(5, 8).Value = 4.5
$ws.Cells.Item(5, 9).Value = 6.5
$ws.Cells.Item(5, 10).Value = 1.95
$ws.Cells.Item(5, 12).Value = 5.5
$ws.Cells.Item(5, 21).Value = 1.67
$ws.Cells.Item(5, 22).Value = 2.1
$ws.Cells.Item(5, 24).Value = 8.5
$ws.Cells.Item(5, 26).Value = 11
$ws.Cells.Item(5, 45).Value = 81
$ws.Cells.Item(5, 48).Value = 41
$ws.Cells.Item(5, 54).Value = 151

$ws.Cells.Item(6, 7).Value = 2.8
$ws.Cells.Item(6, 9).Value = 2.38
$ws.Cells.Item(6, 10).Value = 3.25
$ws.Cells.Item(6, 12).Value = 3
$ws.Cells.Item(6, 24).Value = 15
$ws.Cells.Item(6, 28).Value = 26
$ws.Cells.Item(6, 34).Value = 10
$ws.Cells.Item(6, 35).Value = 13
$ws.Cells.Item(6, 36).Value = 9.5
$ws.Cells.Item(6, 37).Value = 23
$ws.Cells.Item(6, 44).Value = 51
$ws.Cells.Item(6, 50).Value = 13

$ws.Cells.Item(9, 7).Value = 1.33
$ws.Cells.Item(9, 9).Value = 7.5
$ws.Cells.Item(9, 17).Value = 1.3
$ws.Cells.Item(9, 18).Value = 3.5
$ws.Cells.Item(9, 30).Value = 12
$ws.Cells.Item(9, 42).Value = 12
$ws.Cells.Item(9, 44).Value = 26
$ws.Cells.Item(9, 47).Value = 8
$ws.Cells.Item(9, 50).Value = 34
$ws.Cells.Item(9, 51).Value = 29

$ws.Cells.Item(10, 8).Value = 4.1
$ws.Cells.Item(10, 9).Value = 4.2
$ws.Cells.Item(10, 16).Value = 4.4
$ws.Cells.Item(10, 17).Value = 1.53
$ws.Cells.Item(10, 18).Value = 2.35
$ws.Cells.Item(10, 21).Value = 1.6
$ws.Cells.Item(10, 22).Value = 2.22
$ws.Cells.Item(10, 23).Value = 9.5
$ws.Cells.Item(10, 24).Value = 9.25
$ws.Cells.Item(10, 27).Value = 11.75
$ws.Cells.Item(10, 28).Value = 19.5
$ws.Cells.Item(10, 31).Value = 14
$ws.Cells.Item(10, 32).Value = 50
$ws.Cells.Item(10, 33).Value = 300
$ws.Cells.Item(10, 34).Value = 16.5
$ws.Cells.Item(10, 35).Value = 28
$ws.Cells.Item(10, 39).Value = 35
$ws.Cells.Item(10, 42).Value = 14.5
$ws.Cells.Item(10, 47).Value = 7
$ws.Cells.Item(10, 49).Value = 6.4
$ws.Cells.Item(10, 51).Value = 25

# --- Add new row 11 (New Mexico vs Las Vegas Lights) ---
# Column B looks like a date (dd/mm/yyyy); force text format so Excel
# does not auto-convert it into a date serial number.
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "tSZqiGYq"
$ws.Cells.Item(11, 2).Value = "09/11/2024"
$ws.Cells.Item(11, 3).Value = "23:30"
$ws.Cells.Item(11, 4).Value = "USA - USL CHAMPIONSHIP"
$ws.Cells.Item(11, 5).Value = "New Mexico"
$ws.Cells.Item(11, 6).Value = "Las Vegas Lights"
$ws.Cells.Item(11, 7).Value = 2.37
$ws.Cells.Item(11, 8).Value = 3.3
$ws.Cells.Item(11, 9).Value = 2.65
$ws.Cells.Item(11, 10).Value = 2.9
$ws.Cells.Item(11, 11).Value = 2.2
$ws.Cells.Item(11, 12).Value = 3.25
$ws.Cells.Item(11, 13).Value = 1.05
$ws.Cells.Item(11, 14).Value = 7.8
$ws.Cells.Item(11, 15).Value = 1.25
$ws.Cells.Item(11, 16).Value = 3.55
$ws.Cells.Item(11, 17).Value = 1.75
$ws.Cells.Item(11, 18).Value = 2
$ws.Cells.Item(11, 19).Value = 1.34
$ws.Cells.Item(11, 20).Value = 3
$ws.Cells.Item(11, 21).Value = 1.6
$ws.Cells.Item(11, 22).Value = 2.18
$ws.Cells.Item(11, 23).Value = 9.75
$ws.Cells.Item(11, 24).Value = 13.5
$ws.Cells.Item(11, 25).Value = 9.25
$ws.Cells.Item(11, 26).Value = 26
$ws.Cells.Item(11, 27).Value = 18
$ws.Cells.Item(11, 28).Value = 24
$ws.Cells.Item(11, 29).Value = 7.8
$ws.Cells.Item(11, 30).Value = 6.6
$ws.Cells.Item(11, 31).Value = 12.5
$ws.Cells.Item(11, 32).Value = 45
$ws.Cells.Item(11, 33).Value = 300
$ws.Cells.Item(11, 34).Value = 9.75
$ws.Cells.Item(11, 35).Value = 14.5
$ws.Cells.Item(11, 36).Value = 10
$ws.Cells.Item(11, 37).Value = 32
$ws.Cells.Item(11, 38).Value = 21
$ws.Cells.Item(11, 39).Value = 27
$ws.Cells.Item(11, 40).Value = 4.55
$ws.Cells.Item(11, 41).Value = 12.5
$ws.Cells.Item(11, 42).Value = 17.5
$ws.Cells.Item(11, 43).Value = 45
$ws.Cells.Item(11, 44).Value = 70
$ws.Cells.Item(11, 45).Value = 175
$ws.Cells.Item(11, 46).Value = 3
$ws.Cells.Item(11, 47).Value = 6.6
$ws.Cells.Item(11, 48).Value = 50
$ws.Cells.Item(11, 49).Value = 4.8
$ws.Cells.Item(11, 50).Value = 14.5
$ws.Cells.Item(11, 51).Value = 20
$ws.Cells.Item(11, 52).Value = 60
$ws.Cells.Item(11, 53).Value = 90
$ws.Cells.Item(11, 54).Value = 200
$ws.Cells.Item(11, 55).Value = 51
$ws.Cells.Item(11, 56).Value = 51

# Drop the temporary text format on B11 so it matches the plain
# (unstyled) inline-string cells used throughout the rest of the sheet.
$ws.Cells.Item(11, 2).ClearFormats()

